$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.372.35"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "2.584.97"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'552.89"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").Value = "'139.99"
$ws.Range("E6").Value = "  -1.52%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "2.601.51"
$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("E12").Value = "  +5.29%  "

$ws.Range("E13").Value = "  +5.04%  "

$ws.Range("D14").Value = "3.042.46"
$ws.Range("E14").Value = "  +0.43%  "

$ws.Range("D15").Value = "59.368.05"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").Value = "'22.91"
$ws.Range("E16").Value = "  +5.39%  "

$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").Value = "2.592.86"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D20").Value = "'340.17"

$ws.Range("D21").Value = "'10.41"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("D22").Value = "'6.59"
$ws.Range("E22").Value = "  +6.95%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  +8.21%  "

$ws.Range("D25").Value = "'63.01"
$ws.Range("E25").Value = "  -2.57%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("E28").Value = "  +4.22%  "

$ws.Range("E29").Value = "  -1.29%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("D32").Value = "'6.09"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("D33").Value = "'157.46"
$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").Value = "'19.33"
$ws.Range("E34").Value = "  +2.66%  "

$ws.Range("D35").Value = "'4.08"
$ws.Range("E35").Value = "  +2.02%  "

$ws.Range("D36").Value = "'0.911"
$ws.Range("E36").Value = "  +3.73%  "

$ws.Range("E37").Value = "  +3.18%  "

$ws.Range("E38").Value = "  +2.33%  "

$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'0.839"
$ws.Range("E40").Value = "  -4.28%  "

$ws.Range("E41").Value = "  +1.52%  "

$ws.Range("D42").Value = "'288.93"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").Value = "'135.71"
$ws.Range("E43").Value = "  +8.81%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "'0.0973"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "'0.598"
$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("D50").Value = "1.971.45"
$ws.Range("E50").Value = "  +2.88%  "

$ws.Range("E51").Value = "  +1.40%  "
